$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = 0.4643152847966753
$ws.Range("J19").Value = 0.2166321763746371
$ws.Range("K19").Value = 0.09876596577347316
$ws.Range("L19").Value = 2.670391752303761

$ws.Range("I20").Value = 0.7433026278503435
$ws.Range("J20").Value = 0.490140279905831
$ws.Range("K20").Value = 0.4156417725326108
$ws.Range("L20").Value = 2.21484215664195
